# Implement CreateProcess and DeleteProcess workflows
# - Adds new result-column letters to the "Settings" table (Table1)
# - Adds new localized (EN/JA) error-message rows to the "LocalizedSettings"
#   table (Table13) for the new Process/Package validation messages

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
$ws2 = $wb.Worksheets.Item("LocalizedSettings")

$lo1 = $ws1.ListObjects.Item(1)
$lo2 = $ws2.ListObjects.Item(1)

# Grow table2 (LocalizedSettings) by 8 rows (93..100), and table1 (Settings) by 3 rows (37..39)
for ($i = 0; $i -lt 8; $i++) { [void]$lo2.ListRows.Add() }
for ($i = 0; $i -lt 3; $i++) { [void]$lo1.ListRows.Add() }

# --- Fill cells in the same chronological order the shared strings were created ---

# Row 94: ProcessNameNotSpecified block
$ws2.Range("A94").Value = "ProcessNameNotSpecified"
$ws2.Range("B94").Value = "Process name not specified."
$ws2.Range("B94").WrapText = $true
$ws2.Range("C94").Value = "プロセス名が指定されませんでした。"
$ws2.Range("C94").WrapText = $true

# Row 99: PackageNameNotSpecified block
$ws2.Range("A99").Value = "PackageNameNotSpecified"
$ws2.Range("B99").Value = "Package name not specified."
$ws2.Range("B99").WrapText = $true
$ws2.Range("C99").Value = "パッケージ名が指定されませんでした。"
$ws2.Range("C99").WrapText = $true

# Sheet1 row 38: CreateProcessResultColumn
$ws1.Range("A38").Value = "CreateProcessResultColumn"
$ws1.Range("B38").Value = "G"
$ws1.Range("C38").WrapText = $true

# Sheet1 row 37: CreateProcessIDColumn
$ws1.Range("A37").Value = "CreateProcessIDColumn"
$ws1.Range("B37").Value = "F"
$ws1.Range("C37").WrapText = $true

# Row 100: PackageVersionNotSpecified block
$ws2.Range("A100").Value = "PackageVersionNotSpecified"
$ws2.Range("B100").Value = "Package version not specified"
$ws2.Range("B100").WrapText = $true
$ws2.Range("C100").Value = "パッケージバージョンが指定されませんでした。"
$ws2.Range("C100").WrapText = $true

# Row 95: ProcessIDInvalidOrNotSpecified block
$ws2.Range("A95").Value = "ProcessIDInvalidOrNotSpecified"
$ws2.Range("B95").Value = "Process ID invalid or not specified."
$ws2.Range("B95").WrapText = $true
$ws2.Range("C95").Value = "プロセスIDが無効か指定されませんでした。"
$ws2.Range("C95").WrapText = $true

# Row 96: ProcessNotFound block
$ws2.Range("A96").Value = "ProcessNotFound"
$ws2.Range("B96").Value = "Process not found."
$ws2.Range("B96").WrapText = $true
$ws2.Range("C96").Value = "プロセスが見つかりませんでした。"
$ws2.Range("C96").WrapText = $true

# Row 97: ProcessIDAndNameDoNotMatch block
$ws2.Range("A97").Value = "ProcessIDAndNameDoNotMatch"
$ws2.Range("B97").Value = "The specified Process ID and Process name do not match."
$ws2.Range("B97").WrapText = $true
$ws2.Range("C97").Value = "指定されたプロセス名とプロセスIDが一致しません。"
$ws2.Range("C97").WrapText = $true

# Sheet1 row 39: DeleteProcessResultColumn
$ws1.Range("A39").Value = "DeleteProcessResultColumn"
$ws1.Range("B39").Value = "E"
$ws1.Range("C39").WrapText = $true

# Blank separator rows 93 and 98 in LocalizedSettings retain the wrap-text style only
$ws2.Range("B93").WrapText = $true
$ws2.Range("C93").WrapText = $true
$ws2.Range("B98").WrapText = $true
$ws2.Range("C98").WrapText = $true

# Update view/selection state to match the final state captured in the diff:
# the user scrolled/selected near the new rows on LocalizedSettings, then
# returned to the Settings sheet (which stays the active/selected sheet on save,
# with its own selection unchanged at A2).
$ws2.Activate()
$ws2.Range("A98").Select() | Out-Null
$ws1.Activate()
$ws1.Range("A2").Select() | Out-Null
